# Splits the "Video + Reading for ML Pipelines (for Tuesday)" run into three
# runs: "Review slides on " / "github" / " + Reading for ML Pipelines (for Tuesday)"
# on both slides that contain this bullet (slide 3 and slide 8).

$p = $ppt.ActivePresentation

$oldText = "Video + Reading for ML Pipelines (for Tuesday)"
$newLead = "Review slides on "
$newMid  = "github"
$oldLead = "Video"

function Update-Bullet($slideIndex) {
    $s = $p.Slides.Item($slideIndex)
    for ($i = 1; $i -le $s.Shapes.Count; $i++) {
        $shp = $s.Shapes.Item($i)
        if (-not $shp.HasTextFrame) { continue }
        $tr = $shp.TextFrame.TextRange
        $fullText = $tr.Text
        $pos = $fullText.IndexOf($oldText)
        if ($pos -ge 0) {
            # Step 1: replace "Video" with "Review slides on github" in place.
            # This keeps the same run for the untouched remainder
            # (" + Reading for ML Pipelines (for Tuesday)") and creates
            # a fresh run for the newly-typed text.
            $videoRange = $tr.Characters($pos + 1, $oldLead.Length)
            $videoRange.Text = $newLead + $newMid

            # Step 2: split that freshly-typed run so "github" becomes its
            # own run (mirrors PowerPoint auto-flagging it as a possible
            # spelling error once it was typed as a separate word).
            $tr2 = $shp.TextFrame.TextRange
            $ghPos = $tr2.Text.IndexOf($newMid, $pos)
            $ghRange = $tr2.Characters($ghPos + 1, $newMid.Length)
            $ghRange.Text = $newMid
        }
    }
}

Update-Bullet 3
Update-Bullet 8
